# Update the multiplication-fact table cells to the new values.
# Only the five "content" rows (1, 5, 10, 15, 20) of the single
# table in the document carry text; all other rows are blank
# placeholder rows and are left untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("216×5=", "633×5=", "890×7=", "377×8=", "121×5=")
    5  = @("900×3=", "441×3=", "845×8=", "390×4=", "362×6=")
    10 = @("197×4=", "111×2=", "588×8=", "405×8=", "202×7=")
    15 = @("310×6=", "784×8=", "256×8=", "902×8=", "693×2=")
    20 = @("166×8=", "379×9=", "308×3=", "170×9=", "333×2=")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
